# Update countries & provincias Spain
# Applies the data refresh captured in the commit: updated case numbers for a
# handful of countries, two brand-new rows inserted into the (cases-sorted)
# table, the resulting re-sort of countries with tied/near totals, and the
# "last updated" timestamp bump.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $country, $totales, $nuevos, $activos, $recuperados, $criticos, $muertesHoy, $muertes) {
    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $totales
    $ws.Cells.Item($row, 3).Value = $nuevos
    $ws.Cells.Item($row, 4).Value = $activos
    $ws.Cells.Item($row, 5).Value = $recuperados
    $ws.Cells.Item($row, 6).Value = $criticos
    $ws.Cells.Item($row, 7).Value = $muertesHoy
    $ws.Cells.Item($row, 8).Value = $muertes
}

# Timestamp footer (row 1, column A)
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 04:03"

# Brasil - refreshed totals
Set-Row 12 "Brasil" 115455 740 48221 59296 8318 17 7938

# Corea del Sur - refreshed totals
Set-Row 41 "Corea del Sur" 10806 2 9333 1218 55 1 255

# Bolivia enters the ranked table just below Guinea, pushing
# Islandia/Estonia/Bulgaria/Cuba down one row each.
Set-Row 78 "Bolivia" 1802 121 187 1529 3 4 86
Set-Row 79 "Islandia" 1799 0 1733 56 0 0 10
Set-Row 80 "Estonia" 1711 0 261 1395 6 0 55
Set-Row 81 "Bulgaria" 1704 0 342 1282 37 0 80
Set-Row 82 "Cuba" 1685 0 954 662 8 0 69

# Nueva Zelanda - refreshed totals
Set-Row 84 "Nueva Zelanda" 1488 2 1316 151 2 1 21

# Guatemala enters the ranked table just above Niger, pushing
# Costa Rica/Principado de Andorra/Libano/Mayotte down one row each.
Set-Row 101 "Guatemala" 763 33 79 665 5 0 19
Set-Row 102 "Niger" 763 0 543 182 0 0 38
Set-Row 103 "Costa Rica" 755 0 413 336 5 0 6
Set-Row 104 "Principado de Andorra" 751 0 514 191 16 0 46
Set-Row 105 "Libano" 741 0 206 510 43 0 25
Set-Row 106 "Mayotte" 739 0 352 378 6 0 9

# Islas Caimanes moves ahead of San Martin (Parte Holandesa)
Set-Row 168 "Islas Caimanes" 78 3 30 47 3 0 1
Set-Row 169 "San Martin (Parte Holandesa)" 76 0 44 18 7 0 14

# Curazao / Dominica swap order (tied totals)
Set-Row 198 "Curazao" 16 0 13 2 0 0 1
Set-Row 199 "Dominica" 16 0 14 2 0 0 0

# Seychelles / Montserrat swap order (tied totals)
Set-Row 205 "Seychelles" 11 0 8 3 0 0 0
Set-Row 206 "Montserrat" 11 0 7 3 1 0 1
